# Update "想去人数" (column F) counts on the 展览 (sheet1) and 全部类型 (sheet4)
# worksheets to reflect the refreshed scrape output.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 76
$ws1.Range("F3").Value = 11806
$ws1.Range("F5").Value = 346
$ws1.Range("F7").Value = 11739
$ws1.Range("F8").Value = 485
$ws1.Range("F10").Value = 97
$ws1.Range("F11").Value = 41
$ws1.Range("F12").Value = 1772
$ws1.Range("F13").Value = 5814
$ws1.Range("F14").Value = 121
$ws1.Range("F15").Value = 3528

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 76
$ws4.Range("F5").Value = 11806
$ws4.Range("F7").Value = 346
$ws4.Range("F9").Value = 11739
$ws4.Range("F10").Value = 485
$ws4.Range("F12").Value = 97
$ws4.Range("F13").Value = 41
$ws4.Range("F14").Value = 1772
$ws4.Range("F16").Value = 5814
$ws4.Range("F17").Value = 121
$ws4.Range("F18").Value = 3528
